# Update betting-odds figures on the "Jogos da Semana" sheet to the
# refreshed values captured in the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Bucaramanga vs Alianza)
$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.95
$ws.Range("AG2").Value = 19
$ws.Range("AR2").Value = 41

# Row 3 (Llaneros vs Dep. Cali)
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.95
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5

# Row 8 (Crawley vs Wigan)
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("S8").Value = 2.07
$ws.Range("T8").Value = 1.69

# Row 9 (Mansfield vs Lincoln)
$ws.Range("G9").Value = 2.38
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 3.6
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("T9").Value = 1.63
$ws.Range("AA9").Value = 1.95
$ws.Range("AB9").Value = 1.8
$ws.Range("AD9").Value = 11
$ws.Range("AM9").Value = 351
$ws.Range("AP9").Value = 11

# Row 10 (Shrewsbury vs Huddersfield)
$ws.Range("L10").Value = 2.62
$ws.Range("T10").Value = 1.63

# Row 11 (Stevenage vs Burton)
$ws.Range("T11").Value = 1.67

# Row 12 (Wrexham vs Leyton Orient)
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 3.2
$ws.Range("J12").Value = 2.87
$ws.Range("L12").Value = 3.75
$ws.Range("T12").Value = 1.72
$ws.Range("AA12").Value = 1.8
$ws.Range("AB12").Value = 1.95
$ws.Range("AC12").Value = 7.5
$ws.Range("AD12").Value = 10
$ws.Range("AF12").Value = 21
$ws.Range("AG12").Value = 19
$ws.Range("AN12").Value = 9
$ws.Range("AO12").Value = 15
$ws.Range("AQ12").Value = 34
$ws.Range("AR12").Value = 26
$ws.Range("AS12").Value = 34

# Row 13 (Wycombe vs Bristol Rovers)
$ws.Range("J13").Value = 1.87
$ws.Range("S13").Value = 1.67
$ws.Range("W13").Value = 2.63
$ws.Range("X13").Value = 1.44
